$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new blank rows before the current last row (row 60),
# pushing the existing row 60 data down to row 62.
$ws.Rows.Item(60).Resize(2).Insert()

# Step 2: New row 60 - "Forelle" / "Tercera"
$ws.Range("A60").Value = 1
$ws.Range("B60").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C60").Value = "Arica y Parinacota"
$ws.Range("D60").Value = 45041
$ws.Range("D60").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E60").Value = 15
$ws.Range("F60").Value = "Fruta"
$ws.Range("G60").Value = 100104
$ws.Range("H60").Value = "Frutos de pepita"
$ws.Range("I60").Value = 100104005
$ws.Range("J60").Value = "Pera"
$ws.Range("K60").Value = "Forelle"
$ws.Range("L60").Value = "Tercera"
$ws.Range("M60").Value = 220
$ws.Range("N60").Value = 12000
$ws.Range("O60").Value = 13000
$ws.Range("P60").Value = 12455
$ws.Range("Q60").Value = "$/caja 16 kilos granel"
$ws.Range("R60").Value = "Región de O'Higgins"
$ws.Range("S60").Value = 778
$ws.Range("T60").Value = 16

# Step 3: New row 61 - "Packham's Triumph" / "Tercera"
$ws.Range("A61").Value = 1
$ws.Range("B61").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C61").Value = "Arica y Parinacota"
$ws.Range("D61").Value = 45041
$ws.Range("D61").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E61").Value = 15
$ws.Range("F61").Value = "Fruta"
$ws.Range("G61").Value = 100104
$ws.Range("H61").Value = "Frutos de pepita"
$ws.Range("I61").Value = 100104005
$ws.Range("J61").Value = "Pera"
$ws.Range("K61").Value = "Packham's Triumph"
$ws.Range("L61").Value = "Tercera"
$ws.Range("M61").Value = 300
$ws.Range("N61").Value = 12000
$ws.Range("O61").Value = 13000
$ws.Range("P61").Value = 12500
$ws.Range("Q61").Value = "$/caja 16 kilos granel"
$ws.Range("R61").Value = "Región de O'Higgins"
$ws.Range("S61").Value = 781
$ws.Range("T61").Value = 16
